$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells (and two name/link swaps) to match latest scrape.

$ws.Range("D2").Value = "29.168.09"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.859.99"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'0.7134"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "'240.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.07723"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "'0.3075"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "'24.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "'0.08255"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.875.65"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.220"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.7153"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "'90.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "29.173.50"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "'5.864"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'243.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'0.000007805"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'13.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").Value = "2.108.90"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'7.942"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.1573"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "'162.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'8.895"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "'18.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.495"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.317"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").Value = "'4.357"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'4.087"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "'0.05180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").Value = "'1.904"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").Value = "'0.7270"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D38").Value = "'0.01846"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").Value = "'2.686"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "1.146.71"
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D41").Value = "'0.9020"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'6.085"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "'72.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'101.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "2.005.97"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "'0.5238"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").Value = "'1.760"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").Value = "'9.298"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").Value = "'2.866"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
